$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Size = 14
Write-Host "done"
